$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (within the used data range), shifting old D:K data to F:M
$ws.Range("D7:E102").Insert(-4161)

# Copy number formats/styles from column F (which now holds what used to be column D)
# into the newly inserted (blank) columns D:E so they match the original per-row styling
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new quarter data for columns D (period ending 2018-12-31) and E (period ending 2018-09-30)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2200200
$ws.Range("E8").Value = 1273100
$ws.Range("D9").Value = 1584700
$ws.Range("E9").Value = 552900
$ws.Range("D10").Value = 615500
$ws.Range("E10").Value = 720200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 6100
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 111200
$ws.Range("E15").Value = 113800
$ws.Range("D17").Value = 2038600
$ws.Range("E17").Value = 1218800
$ws.Range("D18").Value = 161600
$ws.Range("E18").Value = 54300
$ws.Range("D20").Value = 10500
$ws.Range("E20").Value = 7700
$ws.Range("D21").Value = 283300
$ws.Range("E21").Value = 175800
$ws.Range("D22").Value = 60200
$ws.Range("E22").Value = 57300
$ws.Range("D23").Value = 111900
$ws.Range("E23").Value = 4700
$ws.Range("D24").Value = 23400
$ws.Range("E24").Value = 6700
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 88500
$ws.Range("E26").Value = -2000
$ws.Range("D27").Value = 64200
$ws.Range("E27").Value = 30200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = -5800
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -10500
$ws.Range("E32").Value = -7700
$ws.Range("D33").Value = 64200
$ws.Range("E33").Value = 24400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 64200
$ws.Range("E35").Value = 24400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 477600
$ws.Range("E41").Value = 452600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1209000
$ws.Range("E43").Value = 765900
$ws.Range("D44").Value = 293700
$ws.Range("E44").Value = 318200
$ws.Range("D45").Value = 261900
$ws.Range("E45").Value = 351400
$ws.Range("D46").Value = 2242200
$ws.Range("E46").Value = 1888100
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 5855100
$ws.Range("E48").Value = 5808200
$ws.Range("D49").Value = 3660000
$ws.Range("E49").Value = 3674000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 611000
$ws.Range("E52").Value = 610600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 12368300
$ws.Range("E54").Value = 11980900
$ws.Range("D57").Value = 753300
$ws.Range("E57").Value = 561800
$ws.Range("D58").Value = 695800
$ws.Range("E58").Value = 443700
$ws.Range("D59").Value = 734000
$ws.Range("E59").Value = 726600
$ws.Range("D60").Value = 2183100
$ws.Range("E60").Value = 1732100
$ws.Range("D61").Value = 4150700
$ws.Range("E61").Value = 4146500
$ws.Range("D62").Value = 1987900
$ws.Range("E62").Value = 2002300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 8698900
$ws.Range("E66").Value = 8299500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2620800
$ws.Range("E72").Value = 2610700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3669400
$ws.Range("E76").Value = 3681400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 64200
$ws.Range("E81").Value = 24400
$ws.Range("D83").Value = 111200
$ws.Range("E83").Value = 113800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 96600
$ws.Range("E89").Value = 158500
$ws.Range("D91").Value = -183300
$ws.Range("E91").Value = -180200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -194000
$ws.Range("E94").Value = -175000
$ws.Range("D96").Value = -45300
$ws.Range("E96").Value = -45100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 134000
$ws.Range("E100").Value = -35900
$ws.Range("D101").Value = -3800
$ws.Range("E101").Value = -1700
$ws.Range("D102").Value = 32800
$ws.Range("E102").Value = -54100
